$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new hours value for Week 1 (F4), which feeds G4 via formula =F4
$ws.Range("F4").Value = 5

# Update the view: scroll so C2 is the top-left visible cell, and move the
# active selection to F5
$ws.Range("F5").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 3
